{"js": "// Load all paragraphs in the document body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- 1. Insert a new \"Meta description\" paragraph right after the H1 title ---\nconst titlePara = paragraphs.items[0];\nconst metaPara = titlePara.insertParagraph(\"\", \"After\");\n// New paragraphs inherit the style of the paragraph they split from\n// (Heading 1 here); the inserted paragraph should use the default style.\nmetaPara.style = \"Normal\";\nawait context.sync();\n\n// Bold \"Meta description\" run.\nconst boldRange = metaPara.insertText(\"Meta description\", \"Start\");\nawait context.sync();\nboldRange.font.bold = true;\n\n// Regular run with the rest of the meta description text.\nconst restRange = metaPara.insertText(\n  \": Discover the features of CyberCatz in this slot game review. Play for free and enjoy the unique graphics and free spins feature for a chance to win big.\",\n  \"End\"\n);\nawait context.sync();\nrestRange.font.bold = false;\nawait context.sync();\n\n// --- 2 & 3. Near the end of the document: drop the duplicated bold title\n//            paragraph, and rewrite the italic paragraph's text into the\n//            new feature-image prompt. ---\nconst trailingParagraphs = body.paragraphs;\ntrailingParagraphs.load(\"items\");\nawait context.sync();\n\nconst items = trailingParagraphs.items;\nconst lastPara = items[items.length - 1];\nconst secondLastPara = items[items.length - 2];\n\n// Remove the paragraph that duplicates \"Play CyberCatz Free: Intergalactic Slot Game Review\".\nsecondLastPara.delete();\nawait context.sync();\n\n// Replace the italic description text with the new image-generation prompt,\n// keeping the paragraph's existing (italic) run formatting intact.\nlastPara.insertText(\n  'Create an eye-catching feature image for \"CyberCatz\" featuring a happy Maya warrior wearing glasses in a cartoon style. The image should include elements that reflect the intergalactic and futuristic theme of the game, such as planets, holograms, and hexagonal positions. The Maya warrior should be in a dynamic pose, as if wielding a weapon or casting a spell, to depict the adventurous and action-packed nature of the game. The colors used in the image should be bold and vibrant, grabbing the attention of potential players scrolling through a list of slot games. Overall, the image should entice viewers to click and explore the world of CyberCatz.',\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1. Insert a new \"Meta description\" paragraph right after the H1 title ---\n$titlePara = $d.Paragraphs.Item(1)\n$titlePara.Range.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs.Item(2)\n# A paragraph created by InsertParagraphAfter inherits the style of the\n# paragraph it follows (Heading 1 here); the new paragraph should use the\n# default body style instead.\n$metaPara.Range.Style = \"Normal\"\n\n# Bold \"Meta description\" run.\n$boldRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start)\n$boldRange.InsertAfter(\"Meta description\")\n$boldRange.Bold = 1\n\n# Regular (non-bold) run with the rest of the meta description text.\n$restRange = $d.Range($boldRange.End, $boldRange.End)\n$restRange.InsertAfter(\": Discover the features of CyberCatz in this slot game review. Play for free and enjoy the unique graphics and free spins feature for a chance to win big.\")\n$restRange.Bold = 0\n\n# --- 2 & 3. Near the end of the document: drop the duplicated bold title\n#            paragraph, and rewrite the italic paragraph's text into the\n#            new feature-image prompt. ---\n$count = $d.Paragraphs.Count\n$secondLastPara = $d.Paragraphs.Item($count - 1)\n\n# Remove the paragraph that duplicates \"Play CyberCatz Free: Intergalactic Slot Game Review\".\n$secondLastPara.Range.Delete()\n\n# Replace the italic description text with the new image-generation prompt,\n# keeping the paragraph's existing (italic) run formatting intact. Trim the\n# trailing paragraph-mark character off the range so only the run text (not\n# the mark) is replaced.\n$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$lastParaRange = $lastPara.Range\n$textRange = $d.Range($lastParaRange.Start, $lastParaRange.End - 1)\n$textRange.Text = 'Create an eye-catching feature image for \"CyberCatz\" featuring a happy Maya warrior wearing glasses in a cartoon style. The image should include elements that reflect the intergalactic and futuristic theme of the game, such as planets, holograms, and hexagonal positions. The Maya warrior should be in a dynamic pose, as if wielding a weapon or casting a spell, to depict the adventurous and action-packed nature of the game. The colors used in the image should be bold and vibrant, grabbing the attention of potential players scrolling through a list of slot games. Overall, the image should entice viewers to click and explore the world of CyberCatz.'\n"}
